$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: add D28 (end time) and E28 (duration formula) ---
# Copy number formats from the existing row 27 cells so the new cells
# reuse the same style indices instead of creating new ones.
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 0.75

$ws.Range("E27").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Formula = "=D28-C28"

# --- Row 30 (row 29 intentionally left blank): new date/start entry ---
$ws.Range("B27").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$ws.Range("B30").Value = 42928

$ws.Range("C28").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null
$ws.Range("C30").Value = 0.43055555555555558

$ws.Application.CutCopyMode = $false

$ws.Range("D30").Select()
